$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "w" / "h" columns (E1/F1) for the report-card table ---
$ws.Range("E1").Value = "w"
$ws.Range("F1").Value = "h"

# --- Data rows: widths (E) and heights (F) for the rendered images ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1.5

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1.5

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1

# --- Formatting: make the header row wrap text and taller, matching the rest of the header cells ---
$ws.Range("C1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 33

# --- Selection / view state ---
$ws.Range("F4").Select()
